# Insert a new data row at row 347 (pushing the existing rows 347:479 down
# to 348:480) and populate it with a new "Acelga" price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(347).Insert()

$ws.Cells.Item(347, 1).Value = 10
$ws.Cells.Item(347, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(347, 3).Value = "La Araucanía"
$ws.Cells.Item(347, 4).Value = 45009
$ws.Cells.Item(347, 5).Value = 9
$ws.Cells.Item(347, 6).Value = 100112009
$ws.Cells.Item(347, 7).Value = "Acelga"
$ws.Cells.Item(347, 8).Value = "Sin especificar"
$ws.Cells.Item(347, 9).Value = "Primera"
$ws.Cells.Item(347, 10).Value = 65
$ws.Cells.Item(347, 11).Value = 8000
$ws.Cells.Item(347, 12).Value = 8000
$ws.Cells.Item(347, 13).Value = 8000
$ws.Cells.Item(347, 14).Value = '$/docena de atados (12 kilos)'
$ws.Cells.Item(347, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(347, 16).Value = 667
$ws.Cells.Item(347, 17).Value = 12
$ws.Cells.Item(347, 18).Value = "Hortaliza"
